# Applies the "Updated cryptos list" refresh to the cryptos table on the
# active worksheet: per-row Price (D) / Volume(1h) (E) updates, and a few
# rows where Coin name + Link + Price + Volume were swapped with the
# neighboring row (rows 12/13, 32/33, 46/47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.836.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.027.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.87%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("E9").Value = "  -0.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

# Row 11
$ws.Range("E11").Value = "  +0.84%  "

# Row 12
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.331.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.758"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "

# Row 16
$ws.Range("E16").Value = "  -1.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.030.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.770.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.87%  "

# Row 28
$ws.Range("E28").Value = "  -2.66%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "

# Row 30
$ws.Range("E30").Value = "  -4.34%  "

# Row 31
$ws.Range("E31").Value = "  +0.86%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.27%  "

# Row 33
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0602"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.28%  "

# Row 37
$ws.Range("E37").Value = "  -3.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.519.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "

# Row 41
$ws.Range("E41").Value = "  +0.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$ws.Range("E44").Value = "  -0.29%  "

# Row 45
$ws.Range("E45").Value = "  -1.50%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.55%  "

# Row 47
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.03%  "

# Row 48
$ws.Range("E48").Value = "  -0.76%  "

# Row 49
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.218.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.92%  "
